$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Test Steps" -> add a "Results" column (F) with "PASS" for every
# existing data row.
# ---------------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("Test Steps")

# Header cell F1: copy the look of the existing "highlighted" header style
# (used by A9:A15) and then fill in the text.
$wsSteps.Range("A9").Copy()
$wsSteps.Range("F1").PasteSpecial(-4122)
$wsSteps.Range("F1").Value = "Results"

# Data cells F2:F15: copy the plain bordered look used by the rest of the
# table (e.g. D1) and stamp "PASS" into every row.
$wsSteps.Range("D1").Copy()
$wsSteps.Range("F2:F15").PasteSpecial(-4122)
$wsSteps.Range("F2:F15").Value = "PASS"

$wsSteps.Range("F2").Select()
$wsSteps.Range("F2:F15").Select()

# ---------------------------------------------------------------------------
# Sheet "Test Cases" -> add a matching "Results" column (D) with "PASS" for
# every existing data row.
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("Test Cases")

$wsCases.Range("A1").Copy()
$wsCases.Range("D1").PasteSpecial(-4122)
$wsCases.Range("D1").Value = "Results"

$wsCases.Range("C1").Copy()
$wsCases.Range("D2:D3").PasteSpecial(-4122)
$wsCases.Range("D2:D3").Value = "PASS"

$wsCases.Range("C14").Select()
